$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.414.20"
$ws.Cells.Item(2, 5).Value = "  -0.47%  "

$ws.Cells.Item(3, 4).Value = "1.722.41"
$ws.Cells.Item(3, 5).Value = "  -0.41%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.000"
$ws.Cells.Item(4, 5).Value = "  -0.01%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "243.14"
$ws.Cells.Item(5, 5).Value = "  -1.25%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.001"
$ws.Cells.Item(6, 5).Value = "  +0.03%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4869"
$ws.Cells.Item(7, 5).Value = "  +0.78%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2598"
$ws.Cells.Item(8, 5).Value = "  -2.54%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06199"
$ws.Cells.Item(9, 5).Value = "  -0.63%  "

$ws.Cells.Item(10, 4).Value = "1.726.48"
$ws.Cells.Item(10, 5).Value = "  -0.13%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.06986"
$ws.Cells.Item(11, 5).Value = "  -1.40%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "15.45"
$ws.Cells.Item(12, 5).Value = "  -0.92%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "4.541"
$ws.Cells.Item(13, 5).Value = "  -0.26%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.5975"
$ws.Cells.Item(14, 5).Value = "  -1.81%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "77.42"
$ws.Cells.Item(15, 5).Value = "  +0.23%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "1.0000"
$ws.Cells.Item(16, 5).Value = "  -0.07%  "

$ws.Cells.Item(17, 4).Value = "26.421.72"
$ws.Cells.Item(17, 5).Value = "  -0.43%  "

$ws.Cells.Item(18, 2).Value = "BinanceUSD"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "1.001"
$ws.Cells.Item(18, 5).Value = "  -0.03%  "

$ws.Cells.Item(19, 2).Value = "ShibaInu"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.000007271"
$ws.Cells.Item(19, 5).Value = "  +0.20%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "11.32"
$ws.Cells.Item(20, 5).Value = "  -1.63%  "

$ws.Cells.Item(21, 4).Value = "1.933.96"
$ws.Cells.Item(21, 5).Value = "  -1.13%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.462"
$ws.Cells.Item(22, 5).Value = "  -0.97%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "8.513"
$ws.Cells.Item(23, 5).Value = "  -2.90%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "5.123"
$ws.Cells.Item(24, 5).Value = "  -2.25%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "136.99"
$ws.Cells.Item(25, 5).Value = "  -0.18%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "15.34"
$ws.Cells.Item(26, 5).Value = "  -0.42%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "1.396"
$ws.Cells.Item(27, 5).Value = "  -0.64%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "107.43"
$ws.Cells.Item(28, 5).Value = "  -0.79%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.727"
$ws.Cells.Item(29, 5).Value = "  -2.70%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "3.954"
$ws.Cells.Item(30, 5).Value = "  -0.13%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.07971"
$ws.Cells.Item(31, 5).Value = "  -0.49%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.674"
$ws.Cells.Item(32, 5).Value = "  -0.41%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.04503"
$ws.Cells.Item(33, 5).Value = "  -1.47%  "

$ws.Cells.Item(34, 5).Value = "  -0.64%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.002"
$ws.Cells.Item(35, 5).Value = "  +0.14%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.6285"
$ws.Cells.Item(36, 5).Value = "  -0.46%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.9314"
$ws.Cells.Item(37, 5).Value = "  +4.43%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "1.964"
$ws.Cells.Item(38, 5).Value = "  -1.68%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.392"
$ws.Cells.Item(39, 5).Value = "  -0.51%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.9999"
$ws.Cells.Item(40, 5).Value = "  -0.17%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.01476"
$ws.Cells.Item(41, 5).Value = "  -1.58%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "99.68"
$ws.Cells.Item(42, 5).Value = "  -1.97%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "5.337"
$ws.Cells.Item(43, 5).Value = "  -1.34%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.3839"
$ws.Cells.Item(44, 5).Value = "  -1.35%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "6.852"
$ws.Cells.Item(45, 5).Value = "  -2.67%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.1170"
$ws.Cells.Item(46, 5).Value = "  -0.96%  "

$ws.Cells.Item(47, 5).Value = "  -0.67%  "

$ws.Cells.Item(48, 2).Value = "EnergySwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "7.756"
$ws.Cells.Item(48, 5).Value = "  -1.66%  "

$ws.Cells.Item(49, 2).Value = "Elrond"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "30.24"
$ws.Cells.Item(49, 5).Value = "  -1.12%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.231"
$ws.Cells.Item(50, 5).Value = "  -1.71%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "50.87"
$ws.Cells.Item(51, 5).Value = "  -1.03%  "
